$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 27

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"
$ws.Cells.Item($row, 4).Value = 44610
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = 100114007
$ws.Cells.Item($row, 7).Value = "Jengibre"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 50
$ws.Cells.Item($row, 11).Value = 17000
$ws.Cells.Item($row, 12).Value = 18000
$ws.Cells.Item($row, 13).Value = 17400
$ws.Cells.Item($row, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item($row, 15).Value = "Perú"
$ws.Cells.Item($row, 16).Value = 1338
$ws.Cells.Item($row, 17).Value = 13
$ws.Cells.Item($row, 18).Value = "Hortaliza"
